$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Add-PlainParagraph($text) {
    $endPos = $d.Content.End
    $r = $d.Range($endPos, $endPos)
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $xml = '<w:p ' + $wNs + '><w:r><w:t>' + $escaped + '</w:t></w:r></w:p>'
    $r.InsertXML($xml) | Out-Null
}

function Add-EmptyParagraph() {
    $endPos = $d.Content.End
    $r = $d.Range($endPos, $endPos)
    $xml = '<w:p ' + $wNs + '/>'
    $r.InsertXML($xml) | Out-Null
}

Add-PlainParagraph "Báo cáo tuần 3: test lần 1"
Add-PlainParagraph "Báo cáo tuần 4: test lần 2 để hiểu sâu"
Add-EmptyParagraph

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
